# Update the CDC grant-history table text to cover FY 2012-2016 instead of
# FY 2011-2016 (the 2011 figures were dropped from the underlying data, so
# the descriptive/alt-text strings above the table need to match).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDCgrantHistTable")

# Long descriptive paragraph directly above the table (merged A3:D3).
$ws.Range("A3").Value = "This table shows the grant awards and award dollars CDC made for FY 2012-2016. It is provided as a text alternative to the interactive chart on the CDC page of this website."

# Short "Grant Awards and Award Dollars Description" caption (merged A7:C7).
$ws.Range("A7").Value = "Grant awards and award dollars CDC made for FY 2012-2016."
